$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.402.66'
$ws.Range('E2').Value = '  +2.74%  '
$ws.Range('D3').Value = '2.310.95'
$ws.Range('E3').Value = '  +1.74%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E5').Value = '  +0.69%  '
$ws.Range('D6').Value = '103.96'
$ws.Range('E6').Value = '  +7.06%  '
$ws.Range('E7').Value = '  +1.03%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.531'
$ws.Range('E9').Value = '  +8.45%  '
$ws.Range('D10').Value = '36.75'
$ws.Range('E10').Value = '  +4.60%  '
$ws.Range('D11').Value = '52.80'
$ws.Range('E11').Value = '  +1.18%  '
$ws.Range('D12').Value = '0.0813'
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('E13').Value = '  -1.16%  '
$ws.Range('D14').Value = '7.00'
$ws.Range('E14').Value = '  +2.53%  '
$ws.Range('D15').Value = '2.666.64'
$ws.Range('E15').Value = '  +1.69%  '
$ws.Range('D16').Value = '15.10'
$ws.Range('E16').Value = '  +3.00%  '
$ws.Range('D17').Value = '2.309.83'
$ws.Range('E17').Value = '  +1.85%  '
$ws.Range('D18').Value = '0.810'
$ws.Range('D19').Value = '43.303.17'
$ws.Range('E19').Value = '  +2.84%  '
$ws.Range('E20').Value = '  -0.60%  '
$ws.Range('E21').Value = '  +2.35%  '
$ws.Range('D22').Value = '6.18'
$ws.Range('E22').Value = '  +3.37%  '
$ws.Range('D23').Value = '68.11'
$ws.Range('D24').Value = '242.80'
$ws.Range('E24').Value = '  +2.69%  '
$ws.Range('D25').Value = '2.03'
$ws.Range('E25').Value = '  +2.81%  '
$ws.Range('D26').Value = '2.62'
$ws.Range('E26').Value = '  +1.18%  '
$ws.Range('E27').Value = '  +0.19%  '
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').Value = '3.98'
$ws.Range('E28').Value = '  -0.86%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '24.92'
$ws.Range('E29').Value = '  +5.82%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '2.30'
$ws.Range('E30').Value = '  +8.10%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = '37.07'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').Value = '9.65'
$ws.Range('E32').Value = '  +0.91%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').Value = '167.46'
$ws.Range('E33').Value = '  +2.26%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '5.28'
$ws.Range('E34').Value = '  +0.61%  '
$ws.Range('D35').Value = '18.43'
$ws.Range('E35').Value = '  +4.42%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').Value = '2.53'
$ws.Range('E37').Value = '  +6.80%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.0743'
$ws.Range('E38').Value = '  +1.26%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').Value = '3.06'
$ws.Range('E39').Value = '  -1.27%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '1.88'
$ws.Range('E40').Value = '  +3.33%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.106'
$ws.Range('E41').Value = '  +2.10%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '4.45'
$ws.Range('E42').Value = '  +6.29%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').Value = '0.116'
$ws.Range('E43').Value = '  +0.82%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').Value = '2.70'
$ws.Range('E44').Value = '  +18.88%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '0.0293'
$ws.Range('E45').Value = '  +3.90%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '1.988.49'
$ws.Range('E46').Value = '  +2.02%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '19.05'
$ws.Range('E47').Value = '  +1.43%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '3.06'
$ws.Range('E48').Value = '  +3.41%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '10.01'
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').Value = '55.76'
$ws.Range('E50').Value = '  +3.79%  '
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').Value = '2.95'
$ws.Range('E51').Value = '  +1.61%  '
